# Slide 12: rework the "Pare-feu" area.
#  - Resize/move the SSH picture ("Picture 2") to take the space previously
#    shared with the firewall diagram.
#  - Remove the firewall textbox, firewall diagram picture, the small
#    "Image 10/15/16/18" icons and the straight connectors between them.
#  - Clear the (now orphaned) entrance animations tied to those shapes so
#    the leftover <p:timing> block is dropped entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# Drop all animation effects on this slide first (their targets are the
# shapes we are about to delete, and PowerPoint removes the <p:timing>
# wrapper once the main sequence is empty).
while ($s.TimeLine.MainSequence.Count -gt 0) {
    $s.TimeLine.MainSequence.Item(1).Delete()
}

# Reposition / resize the SSH screenshot picture.
$pic = $s.Shapes.Item("Picture 2")
$pic.Left = 388.74427
$pic.Top = 180.20546
$pic.Width = 340.25576
$pic.Height = 198.48261

# Remove the firewall textbox, firewall picture, the network-icon pictures
# and the connectors joining them.
$namesToDelete = @(
    "ZoneTexte 11",
    "Picture 6",
    "Image 10",
    "Image 15",
    "Image 16",
    "Image 18",
    "Connecteur droit 19",
    "Connecteur droit 23",
    "Connecteur droit 26"
)
foreach ($name in $namesToDelete) {
    $s.Shapes.Item($name).Delete()
}
